$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text in the source data (values such as
# "57.584.34" are not valid numbers). Force every Price cell that we are about
# to rewrite to Text format first so Excel keeps the new value as a literal
# string (preserving things like leading/trailing zeros) instead of coercing
# it into a number.
$textCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D9", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '56.824.74'
$ws.Range('E2').Value = '  +3.58%  '
$ws.Range('D3').Value = '2.365.26'
$ws.Range('E3').Value = '  +3.14%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.54%  '
$ws.Range('D5').Value = '521.21'
$ws.Range('E5').Value = '  +3.29%  '
$ws.Range('D6').Value = '135.06'
$ws.Range('E6').Value = '  +3.37%  '
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +1.90%  '
$ws.Range('D9').Value = '2.362.12'
$ws.Range('E9').Value = '  +2.06%  '
$ws.Range('E10').Value = '  +7.30%  '
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('D12').Value = '5.21'
$ws.Range('E12').Value = '  +6.21%  '
$ws.Range('D13').Value = '0.344'
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('D14').Value = '23.91'
$ws.Range('E14').Value = '  +2.29%  '
$ws.Range('D15').Value = '2.754.10'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('D16').Value = '56.983.66'
$ws.Range('E16').Value = '  +4.03%  '
$ws.Range('E17').Value = '  +3.76%  '
$ws.Range('D18').Value = '2.353.04'
$ws.Range('E18').Value = '  +3.03%  '
$ws.Range('D19').Value = '10.61'
$ws.Range('E19').Value = '  +2.21%  '
$ws.Range('D20').Value = '4.30'
$ws.Range('E20').Value = '  +2.83%  '
$ws.Range('D21').Value = '322.66'
$ws.Range('E21').Value = '  +4.82%  '
$ws.Range('D22').Value = '6.70'
$ws.Range('E22').Value = '  +5.52%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').Value = '61.13'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('E25').Value = '  +6.96%  '
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').Value = '7.82'
$ws.Range('E27').Value = '  +4.81%  '
$ws.Range('D28').Value = '172.33'
$ws.Range('E28').Value = '  -1.58%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').Value = '1.21'
$ws.Range('E29').Value = '  +9.85%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0743'
$ws.Range('E30').Value = '  +3.12%  '
$ws.Range('D31').Value = '6.32'
$ws.Range('E31').Value = '  +4.04%  '
$ws.Range('D32').Value = '1.69'
$ws.Range('E32').Value = '  +4.05%  '
$ws.Range('D33').Value = '18.46'
$ws.Range('E33').Value = '  +3.02%  '
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').Value = '0.961'
$ws.Range('E35').Value = '  +1.25%  '
$ws.Range('D36').Value = '0.994'
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('E37').Value = '  +4.15%  '
$ws.Range('D38').Value = '4.06'
$ws.Range('E38').Value = '  +7.24%  '
$ws.Range('E39').Value = '  +6.80%  '
$ws.Range('D40').Value = '37.60'
$ws.Range('E40').Value = '  +3.67%  '
$ws.Range('E41').Value = '  +1.73%  '
$ws.Range('D42').Value = '139.83'
$ws.Range('E42').Value = '  +10.76%  '
$ws.Range('D43').Value = '3.61'
$ws.Range('E43').Value = '  +5.66%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').Value = '280.05'
$ws.Range('E44').Value = '  +11.02%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '5.21'
$ws.Range('E45').Value = '  +1.47%  '
$ws.Range('D46').Value = '0.0512'
$ws.Range('E46').Value = '  +3.05%  '
$ws.Range('D47').Value = '0.0932'
$ws.Range('E47').Value = '  +3.42%  '
$ws.Range('D48').Value = '0.565'
$ws.Range('E48').Value = '  +2.02%  '
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('E50').Value = '  +4.12%  '
$ws.Range('D51').Value = '17.06'
$ws.Range('E51').Value = '  +2.60%  '
